$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 (shifts all existing rows 4..24 down by one,
# and extends the existing merged cells automatically).
$ws.Rows.Item(4).Insert()

# Copy the formatting (styles) from the row right below (the row that used
# to be row 4, now row 5) onto the newly inserted blank row 4, so the new
# row matches the look of the rest of the data rows.
$ws.Range("A5:N5").Copy()
$ws.Range("A4:N4").PasteSpecial(-4122)
$ws.Rows.Item(4).RowHeight = 24.75

# Re-create the merged cell groups for the new row 4 (matching the pattern
# used by every other data row: B:G, H:K, L:M).
$ws.Range("B4:G4").Merge()
$ws.Range("H4:K4").Merge()
$ws.Range("L4:M4").Merge()

# Populate the new row with the new item's data.
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "ANTINAL 220MG/5ML 60ML SUSP."
$ws.Range("H4").Value = "4:0"
$ws.Range("L4").Value = 24
$ws.Range("N4").Value = "1:0"

# Renumber the "م" sequence column for the rows that were pushed down -
# they each need to be one higher than before the insert (rows 5..23 hold
# what used to be rows 4..22, sequence numbers 1..19, which must become 2..20).
for ($r = 5; $r -le 23; $r++) {
    $ws.Range("A$r").Value = $r - 3
}

# Update the grand-total cell (now on row 24) to include the new row's amount.
$ws.Range("K24").Value = 1470.9500000000001
